$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style of the other headers (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the "Save" values for rows 2-4 (all zero per the diff)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
